$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.626.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.704.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.85%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.69%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.67%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3757"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.74%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3452"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.185"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07456"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.008"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.234"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.936"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.709.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001119"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.005"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06722"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "83.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.365"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.77%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.656.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.452"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.765"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.82%  "

$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.903.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.54%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "130.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.88%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.170"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +18.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.770"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.239"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08784"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.773"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.573"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06517"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02395"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.933"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2195"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.273"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6401"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.003"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6084"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.826"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.127"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.99%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07224"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.87%  "
